$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3249
$ws.Range("F8").Value = 7600
$ws.Range("F13").Value = 138
$ws.Range("F14").Value = 646
$ws.Range("F15").Value = 1096
$ws.Range("F18").Value = 162
$ws.Range("F19").Value = 1391
$ws.Range("G19").Value = 138
$ws.Range("F21").Value = 6014
$ws.Range("F22").Value = 22
$ws.Range("F24").Value = 4156
$ws.Range("F25").Value = 2991
$ws.Range("F26").Value = 270
$ws.Range("F27").Value = 82
$ws.Range("F28").Value = 82
$ws.Range("F29").Value = 1020
$ws.Range("F33").Value = 1021
$ws.Range("F42").Value = 571
$ws.Range("F43").Value = 364
$ws.Range("F44").Value = 299
$ws.Range("F45").Value = 1050
$ws.Range("F46").Value = 461
$ws.Range("F47").Value = 4
$ws.Range("F48").Value = 1857
$ws.Range("F49").Value = 56
$ws.Range("F50").Value = 313

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 593
$ws.Range("F12").Value = 70
$ws.Range("F15").Value = 169
$ws.Range("F19").Value = 144
$ws.Range("F27").Value = 4436
$ws.Range("F28").Value = 4436
$ws.Range("F29").Value = 4
$ws.Range("F36").Value = 50

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1937
$ws.Range("F8").Value = 3011
$ws.Range("F10").Value = 1243
$ws.Range("F13").Value = 2024
$ws.Range("F14").Value = 8698
$ws.Range("F15").Value = 848

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 3249
$ws.Range("F6").Value = 1937
$ws.Range("F8").Value = 3011
$ws.Range("F9").Value = 7600
$ws.Range("F11").Value = 1243
$ws.Range("F16").Value = 138
$ws.Range("F17").Value = 848
$ws.Range("F19").Value = 593
$ws.Range("F20").Value = 593
$ws.Range("F21").Value = 646
$ws.Range("F22").Value = 1096
$ws.Range("F24").Value = 70
$ws.Range("F25").Value = 162
$ws.Range("F26").Value = 169
$ws.Range("F27").Value = 1391
$ws.Range("G27").Value = 138
$ws.Range("F29").Value = 6014
$ws.Range("F31").Value = 4156
$ws.Range("F32").Value = 3014
$ws.Range("F33").Value = 270
$ws.Range("F34").Value = 82
$ws.Range("F35").Value = 82
$ws.Range("F37").Value = 1021
$ws.Range("F42").Value = 144
$ws.Range("F44").Value = 571
$ws.Range("F45").Value = 299
$ws.Range("F47").Value = 461
$ws.Range("F48").Value = 1857
$ws.Range("F49").Value = 56
$ws.Range("F50").Value = 4436
$ws.Range("F51").Value = 4
